# Applies updated transition-matrix probabilities to Sheet1.
# Values correspond to the recomputed probabilities after adding more games,
# speeding up the simulate-game logic, and drafting optimization logic.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "B2" = 0.2189542483660131
    "C2" = 0.5196078431372549
    "J2" = 0.006535947712418301
    "P2" = 0.1601307189542484
    "S2" = 0.09477124183006536
    "B3" = 0.01234567901234568
    "C3" = 0.01851851851851852
    "J3" = 0.02469135802469136
    "P3" = 0.7901234567901234
    "S3" = 0.154320987654321
    "J4" = 0.01818181818181818
    "P4" = 0.7454545454545455
    "S4" = 0.2363636363636364
    "B6" = 0.04060913705583756
    "D6" = 0.02538071065989848
    "F6" = 0.05583756345177665
    "J6" = 0.3248730964467005
    "O6" = 0.01015228426395939
    "Q6" = 0.16751269035533
    "R6" = 0.04060913705583756
    "S6" = 0.3350253807106599
    "B7" = 0.09293680297397769
    "D7" = 0.02973977695167286
    "E7" = 0.003717472118959108
    "F7" = 0.03345724907063197
    "J7" = 0.1189591078066914
    "O7" = 0.04089219330855019
    "Q7" = 0.1970260223048327
    "R7" = 0.06691449814126393
    "S7" = 0.4163568773234201
    "B8" = 0.09932279909706546
    "D8" = 0.02031602708803612
    "F8" = 0.05417607223476298
    "J8" = 0.1060948081264108
    "O8" = 0.02031602708803612
    "Q8" = 0.1941309255079007
    "R8" = 0.08126410835214447
    "S8" = 0.4243792325056434
    "B9" = 0.1116504854368932
    "D9" = 0.02912621359223301
    "E9" = 0.004854368932038835
    "F9" = 0.04854368932038835
    "J9" = 0.1213592233009709
    "O9" = 0.01941747572815534
    "Q9" = 0.2087378640776699
    "R9" = 0.07281553398058252
    "S9" = 0.383495145631068
    "B10" = 0.1297709923664122
    "D10" = 0.02862595419847328
    "E10" = 0.001908396946564885
    "F10" = 0.06202290076335878
    "J10" = 0.0916030534351145
    "O10" = 0.01717557251908397
    "Q10" = 0.2194656488549618
    "R10" = 0.07538167938931298
    "S10" = 0.3740458015267176
    "G11" = 0.1310160427807487
    "J11" = 0.053475935828877
    "K11" = 0.1711229946524064
    "L11" = 0.6283422459893048
    "S11" = 0.0160427807486631
    "F12" = 0.00392156862745098
    "G12" = 0.7607843137254902
    "J12" = 0.1450980392156863
    "K12" = 0.007843137254901961
    "L12" = 0.05882352941176471
    "S12" = 0.02352941176470588
    "G13" = 0.6363636363636364
    "J13" = 0.3272727272727273
    "S13" = 0.03636363636363636
    "F15" = 0.04017857142857143
    "H15" = 0.1339285714285714
    "I15" = 0.1026785714285714
    "J15" = 0.3303571428571428
    "K15" = 0.04017857142857143
    "M15" = 0.01339285714285714
    "O15" = 0.04910714285714286
    "S15" = 0.2901785714285715
    "F16" = 0.03365384615384615
    "H16" = 0.1490384615384615
    "I16" = 0.07692307692307693
    "J16" = 0.3269230769230769
    "K16" = 0.1538461538461539
    "M16" = 0.02403846153846154
    "O16" = 0.03846153846153846
    "S16" = 0.1971153846153846
    "F17" = 0.0178173719376392
    "H17" = 0.2160356347438753
    "I17" = 0.1158129175946548
    "J17" = 0.3028953229398664
    "K17" = 0.1180400890868597
    "M17" = 0.0311804008908686
    "N17" = 0.0022271714922049
    "O17" = 0.06013363028953229
    "S17" = 0.1358574610244989
    "F18" = 0.03267973856209151
    "H18" = 0.2418300653594771
    "I18" = 0.08496732026143791
    "J18" = 0.3137254901960784
    "K18" = 0.1372549019607843
    "M18" = 0.0261437908496732
    "O18" = 0.07843137254901961
    "S18" = 0.08496732026143791
    "F19" = 0.01721439749608764
    "H19" = 0.1964006259780908
    "I19" = 0.08528951486697965
    "J19" = 0.3028169014084507
    "K19" = 0.1471048513302035
    "M19" = 0.02269170579029734
    "N19" = 0.001564945226917058
    "O19" = 0.07511737089201878
    "S19" = 0.1517996870109546
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
